$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 becomes styled like the rest of the table (style 1)
$ws.Range("A34:I34").Style = "Normal"

# Row 35: new "Extra Trees" result row
$ws.Range("A35").Value = "Extra Trees"
$ws.Range("B35").Value = 0.93
$ws.Range("C35").Value = 0.93
$ws.Range("D35").Value = 0.87
$ws.Range("E35").Value = 0.93
$ws.Range("F35").Value = 0.92
$ws.Range("G35").Value = 0.94
$ws.Range("H35").Value = 0.93
$ws.Range("I35").Value = 0.91

# Rows 36-44: clear the leftover blank styled cells
$ws.Range("A36:I44").Clear()
$ws.Rows("36:44").RowHeight = 13.5

# View state: scroll position + active cell selection
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("G66").Select()
